$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "RM6 PRO LONG TERM REVIEW Having used the Realme 6 pro for 30 days intensly ,I guess I can help all of you out there to make a well informed decison :) .I'm taking this effort so that you’ll know whether your money is worth investing on this phone.Please read on if you consider buying it.1) DESIGN : 4.5/5 .The phone feels premium and I love the flashy design of realme 6 pro.The weight is manageable once you get used to it .The design is subjective and is a matter of personal preference t..."
$ws.Range("B3").Value = "Undoubtedly this is a masterpiece. I took about a year to decide and am COMPLETELY satisfied. Here are the details:* Fantastic and perfect cameras. 100% value for your money. PERFECT results. I am over satisfied (I am also a DSLR photographer). The zoom is REALLY good, try it in the daylight. Night mode is much superior as compared to the other handsets. The best part is that it has a front-facing ultra-wide-angle selfie camera too. Trust me, you'll NEVER regret once you buy it. Everything ..."
$ws.Range("B4").Value = "I received this on my birthday . Real happy. Good camera . Battery life is also very good. It takes 1 hours to full charge . We can play pubg smoothly . Happy to buy from Flipkart . Thanks Flipkart highly recommend to buy this phone.value for money phone ."
$ws.Range("B5").Value = "I've been using Realme 6 Pro (6GB/128GB) since 10 days & from day 1 it had been a heck of an experience. I must say I've never seen such an amazing device in a decade, The uniqueness of this device is exactly what grabs the interest of the user.                      The performance, the outlook ! Wow... Realme just did it. The 64 MP AI Quad Cam is a beast in itself.. and it feels good to have a Qual comm Snapdragon 720g processor I mean it's the next level. And how could I even forget to men..."
$ws.Range("B6").Value = "Good camera .. but not the best as images are sometimes noisy when zoom in for a bit.. selfie camera is better than rear camera.. fast charging and battery backup are really good.. overall a good phone .. but slightly overprice.. here are my best shots from the mobile"
$ws.Range("B7").Value = "Worth it....1stly I Will talk about Flipkart service it is awesome everytime....Mobile is superb...NYC battery backup NYC display smooth touch superb design  ...Camera is also superb but after 20x zoom it doesn't give clear pic ...Mobile don't get heated after heavy use also for this a big thumbs up....Total I think it is a superb mob.... don't think go for it just buy ittttt....Wowwwww"
$ws.Range("B8").Value = "Best in segment todayDisplay is good visibility is nice but this phone should have amoled display for better viewing experience because amoled is better far better than 90hz display but 90hz is super smooth 120hz is battery killer but 90 is perfect for this phone. Camera is best and full of features not as sharp and crisp as poco x2 because of sony sensor. it has tons of features will upload a camera review on YouTube. Portrait blur in video mode is awesome front slow motion and wide ang..."
$ws.Range("B9").Value = "Ok...this is literally the best phone in the segment. Battery is insane, performance/gaming is fantastic, cameras are incredible and overall this phone is just amazing value for money. Although the POCO X2 has a bigger battery than this phone, it's 120hz display sucks alot of power and drains is quicker so battery is better on this phone....."
$ws.Range("B10").Value = "Awesome camera, front and back camera good performance.good battery backup, display so smooth, no lag for pubg, worth for mony,, thank you flipkart,,, love from tamil nadu..."
$ws.Range("B11").Value = "Detailed Review - I have Been Using this phone from last 4 DaysExcellent Delivery By Flipkart Less Than 24HrsPros - 1.  Battery Lasts More Than 46Hrs for Heavy Usage Like Call, FB, Twitter, Live TV Etc2. Performance - 5 Star - Not even a single second Lag3. Camera -4.5 Star - A bit Photo Blurr of 64 MP Camera, Else than that Every Thing Is Perfect. 4. UI - Same Like Apple iPhone Controls are Fascinating. Much Better Than Color OS.Cons - Till Date NO, Will Update Soon"
$ws.Range("B12").Value = "I think it is best for it's price The camera is awesome The SD 720G is very thermal and battery efficient It gives an excellent performance while gaming I think the diaplay is not such great, the colours in the diaplay is not so much accurate like poco x2, but the 90 hetz display is awesome If I talk about the camera quality, it is also best in it's price range. I think realme had not added the camera sensors just to increase the count specially the 12mp telephoto sensor is more acceptab..."
$ws.Range("B13").Value = "It's the best... I'm experiencing both poco x2 and realme 6 pro And I would say both are very good.. If you want a futuristic WELL BALANCED, GOOD QUALITY, DURABLE PHONE with good battery life pay a little bit extra.. go with realme otherwise if you want a good gaming oriented phone and battery doesn't matter for you go with (you know)... Reasons to buy :--New processor always better than its small successer -navic support -future ready 90 hz display (that is greater than 60hz) -Good (..."
$ws.Range("B14").Value = "The Camera📸(64 MP quad cam - Samsung Icocell GW1 sensor)provided is well satisfactory with great stabilisation features ,modes  and camera composition.You could have 20x zoom but don't expect clarity with increasing zoom.It is interesting to have 12 MP telephoto lens which benefits vlogging.There is a slight camera bump (not an issue). When it comes to the front cam ,you have got a dual selfie cam with 16MP primary lens and additional 8MP wide lens ,so u could have greater coverage or field ..."
$ws.Range("B15").Value = "I have been using this phone extensively for a week.1. First the camera which is decent for the price segment. 64 MP clicks nice , clear pictures. Portrait mode works perfect. Macro lens is useless and clicks very bad resolution photos. The zoom is good and thanks to 12 MP telephoto lens for that. Front selfie camera is very good. Ultra wide photos from the frontal camera do struggle in low lighting but they are good when the light is perfect. No noise in the photos from both rear and front ..."
$ws.Range("B16").Value = "🚫Only problem is that it is not communicating with ""Navigation with Indian constellation"" (NavIC) as per the test result by ""GPSTest"".(Edit) 👉Thanks realme - my problem is resolved with the update on 17th March 2020. Really appreciate it. Thanks again.Otherwise:-✅ Value for money- 🌟🌟🌟🌟🌟= 5 stars (It eliminate all of it's competitors by giving some best flagship features (Excellent camera,God like gaming performance,Buttery smooth display due 90 hz refresh rate same as OnePlus 7 pro..."
$ws.Range("B17").Value = "Features are Good, and the camera is decent and the built quality with that lighting design is eye catching,The camera has 20x zoom which does a great job in well lighting conditions, the ram management is doing a good job, and coming to the processor""720G"" is good in managing the graphics with no lags.Battery is really good in this device and it comes all day even if you are a heavy user, It also charges from 0% to 100% in 1hr 10 min. And the realme UI is really done well. Mobile network..."
$ws.Range("B18").Value = "Pros ::1.90hz screen refresh rate was good and 120hz touch sampling rate.2. the gaming experience was good.3.the telephoto lens was an advantage.Cons ::1.Camera quality was not that good even of 64MP Samsung sensor they should have gone for the Sony instead2.Ram management was poor bcoz there is only 2.9GB available out of 6GB and no background apps running3.the realme UI was not that good not in terms of features but in terms of speed & smooth like Oxygen os4.the full screen was n..."
$ws.Range("B19").Value = "The phone is awesome but display is not good , the display panel is poor and there is  a red tone , I noticed some screen bleed outside the front camera , but camera awesome , performance best is this price segment , fast charger is also very good battery backup is also good , looks also very good , Overall I have to say that If you want a good display phone then you don't go with this , if you need better camera and gaming performance then you will go with this phone ."
$ws.Range("B20").Value = "I have looking some feature loaded phones. I have an option of few more phone in my list but I chose realme 6 pro. i an good with this phone. Battery is too good. Performance is fantastic. Display is mind blowing and main features of this phone 90hz Refresh rate display and with good package of cameras. I am loving it this phone of every features. overall best basic flagship phone. you gonna love it this phone. but cons of this phone is heavy phone is some long time holding this phone but it ..."
$ws.Range("B21").Value = "Battery drain so fast as I am a usual gamer battery is not giving good backup.Camera is not to good only day photos are good but indoor photo are not too good.Processor is so powerful and excellent pubg playing experience on this awesome.Display is also cool however it is LCD screen display but colour Chromecast is managed so good.At last I just want so say value for money is not up to the expectation"
